# Regenerate sval data to filter save games.
# Updates the raw per-game stat columns (TB, d2S, K, IP) and the
# computed "sum" column (G = B + C + D + E) for rows 2-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1169995834814548, 9.983522426115931,  3.223369029078222, 13.86384647080068,  27.18773750947629)
    3 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.178645819794754)
    4 = @(0.6545652718822623, 1.626987699542094,  3.223369029078222, 0.5333859586016987, 6.038307959104277)
    5 = @(0.2881169905109251, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 4.349763226824225)
    6 = @(0.2881169905109251, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1.84748871573303)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B: TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C: d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D: K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E: IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G: sum
}
